$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.021.33"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.650.86"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.59"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5217"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2612"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06267"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.43"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.658.96"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "1.879.72"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5414"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "0.0₅8075"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.97"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "26.025.73"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.564"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.15"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.998"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.974"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.81"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.249"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.10"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.405"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05941"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.273"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.477"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.227"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.538"
$ws.Range("E34").Value = "  -7.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.411"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9433"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5777"
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01596"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.851"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8441"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.56"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "1.002.27"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").Value = "1.792.80"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.57"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.924"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4292"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.476"
$ws.Range("E51").Value = "  +0.00%  "
